# Updated cryptos list with refreshed price and volume(1h) figures.
# Source data is textual (prices use "." as thousands separators in
# some rows), so number formats are forced to Text before assignment
# to prevent Excel from auto-converting the strings to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.318.17"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.590.33"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.30"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0610"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.40"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.92"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.57"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.324.18"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.95"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.09"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.329.54"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("E42").Value = "  -24.84%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.765"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.725.78"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.10"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0979"
$ws.Range("E50").Value = "  -4.17%  "
$ws.Range("E51").Value = "  -0.28%  "
